$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column E values to mirror column D entries (attendance counts added)
$ws.Range("E10").Value = 19
$ws.Range("E11").Value = $ws.Range("D11").Value2
$ws.Range("E12").Value = 6
$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("E18").Value = 3
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 3

# Match column D's style/font for the newly-filled cells in rows 17-21
$ws.Range("E17:E21").Font.Size = $ws.Range("D17").Font.Size

# Move the active cell selection to G19 (was G18)
$ws.Range("G19").Select()
